$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove C2 and E2 entirely (naive component forecaster bug fix)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: remove C3 entirely; update E3
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 0.6970543652217387

# Row 4: update C4
$ws.Range("C4").Value = -0.01587181126745385

# Row 5: update C5, E5
$ws.Range("C5").Value = -0.02256889165886955
$ws.Range("E5").Value = -0.02753509623224515

# Row 6: update C6
$ws.Range("C6").Value = 0.09611428386595566

# Row 7: update E7
$ws.Range("E7").Value = -0.2251688766574889

# Row 8: update C8
$ws.Range("C8").Value = -0.001350220946472191

# Row 10: update C10
$ws.Range("C10").Value = -0.5761528471665334

# Row 14: update C14
$ws.Range("C14").Value = -0.4278219446121501

# Row 15: update C15
$ws.Range("C15").Value = -1.026566979837429

# Row 17: update C17, E17
$ws.Range("C17").Value = 0.4636049209196802
$ws.Range("E17").Value = 0.2986939435938973

# Row 18: update C18, E18
$ws.Range("C18").Value = 0.6216390921348403
$ws.Range("E18").Value = -0.0776179936130994

# Row 19: update C19
$ws.Range("C19").Value = -0.6768900623516871
